# Add a "Deny" column (G) to the next-possible-queues Drools decision table.
# Mirrors the existing "Default next queue" / "Default return queue" columns
# (D/E/F) by adding a new action column G wired to
# $model.setDefaultDenyQueue($param) and fixing the F column header, which
# was mislabeled "Default next queue" even though it drives
# setDefaultReturnQueue.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11 (Sequential = true, now spanning through column G) ---
$ws.Range("F11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
# Leading apostrophe forces literal text "true" (matching D11/E11/F11, which
# are text, not boolean); re-paste the format afterwards to drop the
# transient quote-prefix style the text-entry leaves behind.
$ws.Range("G11").Value2 = "'true"
$ws.Range("F11").Copy()
$ws.Range("G11").PasteSpecial(-4122)

# --- Row 12 (blank separator row) : drop the stray bordered cell in G ---
$ws.Range("G12").Clear()

# --- Row 13 (RuleTable title row) ---
$ws.Range("F13").Copy()
$ws.Range("G13").PasteSpecial(-4122)

# --- Row 14 (CONDITION / ACTION header row) ---
$ws.Range("F14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value2 = "ACTION"

# --- Row 15 ($model: NextPossibleQueuesModel row, blank in D/E/F/G) ---
$ws.Range("F15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# --- Row 16 (snippet row: the actual Drools action snippets) ---
$ws.Range("F16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value2 = '$model.setDefaultDenyQueue($param);'

# --- Row 17 (column header labels) ---
# Copy G17's format from F17 *before* relabeling F17, so both end up using
# the same "action header" style already used by D17/E17/F17.
$ws.Range("F17").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value2 = "Default deny queue"
$ws.Range("F17").Value2 = "Default return queue"

# --- Row 18 (the single "Default Workflow" rule row) ---
$ws.Range("F18").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("G18").Value2 = "null"

$excel.CutCopyMode = 0
